$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 28.453125
$ws.Columns.Item(8).ColumnWidth = 106.453125

# --- Update cell values for rows 15-56 (binning backlog items) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'Confirm Input Dataset'
$ws.Range("C15").Value = 'Bin dataset & save to shared storage'
$ws.Range("D15").Value = '-'
$ws.Range("E15").Value = 'Closed'
$ws.Range("F15").Value = 45004
$ws.Range("H15").Value = '-'

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 'Good/Bad Definition'
$ws.Range("C16").Value = 'Show define bad & indeterminate & good definition UI'
$ws.Range("D16").Value = '-'
$ws.Range("E16").Value = 'Closed'
$ws.Range("F16").Value = 45003
$ws.Range("H16").Value = '-'

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 'Good/Bad Definition'
$ws.Range("C17").Value = 'Show confirm definitions button'
$ws.Range("D17").Value = '-'
$ws.Range("E17").Value = 'Closed'
$ws.Range("F17").Value = 45001
$ws.Range("H17").Value = '-'

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 'Good/Bad Definition'
$ws.Range("C18").Value = 'Add definition row when add button is clicked for bad & indeterminate definition'
$ws.Range("D18").Value = '-'
$ws.Range("E18").Value = 'Closed'
$ws.Range("F18").Value = 45003
$ws.Range("H18").Value = '-'

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 'Good/Bad Definition'
$ws.Range("C19").Value = 'Remove definition row when remove button is clicked for bad & indeterminate definitions'
$ws.Range("D19").ClearContents()
$ws.Range("E19").Value = 'Pending'
$ws.Range("F19").ClearContents()
$ws.Range("H19").Value = 'Require installation of dash-renderer in code environment, not sure why but even in local computer, intallation gets error'

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 'Good/Bad Definition'
$ws.Range("C20").Value = 'Save numerical definition when confirm button is clicked'
$ws.Range("D20").Value = '-'
$ws.Range("E20").Value = 'Closed'
$ws.Range("F20").Value = 45003
$ws.Range("H20").Value = '-'

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 'Good/Bad Definition'
$ws.Range("C21").Value = 'Save categorical definition when confirm button is clicked'
$ws.Range("D21").Value = '-'
$ws.Range("E21").Value = 'Closed'
$ws.Range("F21").Value = 45003
$ws.Range("H21").Value = '-'

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 'Good/Bad Definition'
$ws.Range("C22").Value = 'Show only numerical variable in dataset as dropdown options for bad numerical variables'
$ws.Range("D22").Value = '-'
$ws.Range("E22").Value = 'Closed'
$ws.Range("F22").Value = 45001
$ws.Range("H22").Value = '-'

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 'Good/Bad Definition'
$ws.Range("C23").Value = 'Show only categorical variable in dataset as dropdown options for bad categorical variables'
$ws.Range("D23").Value = '-'
$ws.Range("E23").Value = 'Closed'
$ws.Range("F23").Value = 45001
$ws.Range("H23").Value = '-'

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 'Good/Bad Definition'
$ws.Range("C24").Value = 'Update categorical options based on dropdown value'
$ws.Range("D24").Value = '-'
$ws.Range("E24").Value = 'Pending'
$ws.Range("F24").ClearContents()
$ws.Range("H24").Value = 'Require installation of dash-renderer in code environment, not sure why but even in local computer, intallation gets error'

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 'Good/Bad Definition'
$ws.Range("C25").Value = 'Plan data structure for good bad definition'
$ws.Range("D25").Value = '-'
$ws.Range("E25").Value = 'Closed'
$ws.Range("F25").Value = 45003
$ws.Range("H25").Value = '-'

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 'Good/Bad Definition'
$ws.Range("C26").Value = 'Merge numerical ranges before save'
$ws.Range("D26").Value = '-'
$ws.Range("E26").Value = 'Closed'
$ws.Range("F26").Value = 45003
$ws.Range("H26").Value = '-'

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 'Good/Bad Definition'
$ws.Range("C27").Value = 'Merge categorical elements before save'
$ws.Range("D27").Value = '-'
$ws.Range("E27").Value = 'Closed'
$ws.Range("F27").Value = 45003
$ws.Range("H27").Value = '-'

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 'Good/Bad Definition'
$ws.Range("C28").Value = 'Ensure upper bound > lower bound for numerical variables before save, does not save if violated'
$ws.Range("D28").Value = '-'
$ws.Range("E28").Value = 'Closed'
$ws.Range("F28").Value = 45004
$ws.Range("H28").Value = '-'

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 'Good/Bad Definition'
$ws.Range("C29").Value = 'If upper bound > lower bound is violated, show error message'
$ws.Range("D29").Value = '-'
$ws.Range("E29").Value = 'Closed'
$ws.Range("F29").Value = 45004
$ws.Range("H29").ClearContents()

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 'Good/Bad Definition'
$ws.Range("C30").Value = 'Validate any overlapping numerical user input between bad & indeterminate before data is saved'
$ws.Range("D30").Value = '-'
$ws.Range("E30").Value = 'Closed'
$ws.Range("F30").Value = 45004

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 'Good/Bad Definition'
$ws.Range("C31").Value = 'Validate any overlapping categorical user input between bad & indeterminate before data is saved'
$ws.Range("D31").Value = '-'
$ws.Range("E31").Value = 'Closed'
$ws.Range("F31").Value = 45004

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 'Good/Bad Definition'
$ws.Range("C32").Value = 'Show error message when overlapping numerical user input'
$ws.Range("D32").Value = '-'
$ws.Range("E32").Value = 'Closed'
$ws.Range("F32").Value = 45004

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'Good/Bad Definition'
$ws.Range("C33").Value = 'Show error message when overlapping categorical user input'
$ws.Range("D33").Value = '-'
$ws.Range("E33").Value = 'Closed'
$ws.Range("F33").Value = 45004

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 'Good/Bad Definition'
$ws.Range("C34").Value = 'Show statistics with class GoodBadCounter'
$ws.Range("D34").Value = '-'
$ws.Range("E34").Value = 'Closed'
$ws.Range("F34").Value = 45004

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 'Good/Bad Definition'
$ws.Range("C35").Value = 'Show bar chart'
$ws.Range("D35").Value = '-'
$ws.Range("E35").Value = 'Closed'
$ws.Range("F35").Value = 45004

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 'Good/Bad Definition'
$ws.Range("C36").Value = 'Compute sample bad count'

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 'Good/Bad Definition'
$ws.Range("C37").Value = 'Compute sample indeterminate count'

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'Good/Bad Definition'
$ws.Range("C38").Value = 'Compute sample good count, population good count, and population bad count'

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = '(Automated Binning)'
$ws.Range("C39").Value = 'Perform no binning'
$ws.Range("D39").Value = '-'
$ws.Range("E39").Value = 'Closed'
$ws.Range("F39").Value = 45004

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = '(Automated Binning)'
$ws.Range("C40").Value = 'Perform equal-width binning with width'

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = '(Automated Binning)'
$ws.Range("C41").Value = 'Perform equal-width binning with number of bins'

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = '(Automated Binning)'
$ws.Range("C42").Value = 'Perform equal-frequency binning with frequency'

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = '(Automated Binning)'
$ws.Range("C43").Value = 'Perform equal-frequency binning with number of bins'

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = '(Automated Binning)'
$ws.Range("C44").Value = 'Perform binning with bins settings'

$ws.Range("A45").Value = 44

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 'Preview & Download Settings'
$ws.Range("C46").Value = 'Show preview dataset based on binned_df in storage'
$ws.Range("D46").Value = '-'
$ws.Range("E46").Value = 'Closed'
$ws.Range("F46").Value = 45004

$ws.Range("A47").Value = 46

$ws.Range("A48").Value = 47

$ws.Range("A49").Value = 48

$ws.Range("A50").Value = 49

$ws.Range("A51").Value = 50

$ws.Range("A52").Value = 51

$ws.Range("A53").Value = 52

$ws.Range("A54").Value = 53

$ws.Range("A55").Value = 54

$ws.Range("A56").Value = 55
